# Auto-generated Excel COM-interop script
# Applies scheduled price-runner updates to the Brynhildr_Profits leve-profit tables
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1823
$ws.Range("J18").Value = 1201
$ws.Range("L18").Value = 1201
$ws.Range("N18").Value = -1769

$ws.Range("H98").Value = 2137.3333
$ws.Range("I98").Value = 2243.2307
$ws.Range("J98").Value = 1449
$ws.Range("K98").Value = 2243.2307
$ws.Range("L98").Value = 1449
$ws.Range("M98").Value = -745.2307000000001
$ws.Range("N98").Value = -4445

$ws.Range("H107").Value = 4108.6313
$ws.Range("I107").Value = 4209.2144
$ws.Range("J107").Value = 3827
$ws.Range("K107").Value = 4209.2144
$ws.Range("L107").Value = 3827
$ws.Range("M107").Value = -2289.2144
$ws.Range("N107").Value = -7667

$ws.Range("H122").Value = 2137.3333
$ws.Range("I122").Value = 2243.2307
$ws.Range("J122").Value = 1449
$ws.Range("K122").Value = 6729.6921
$ws.Range("L122").Value = 4347
$ws.Range("M122").Value = -4279.6921
$ws.Range("N122").Value = -9247

$ws.Range("H136").Value = 99999
$ws.Range("J136").Value = 99999
$ws.Range("L136").Value = 99999
$ws.Range("N136").Value = -110199

$ws.Range("H137").Value = 5152.3955
$ws.Range("I137").Value = 1856.8276
$ws.Range("J137").Value = 11978.929
$ws.Range("K137").Value = 5570.4828
$ws.Range("L137").Value = 35936.787
$ws.Range("M137").Value = -3020.4828
$ws.Range("N137").Value = -41036.787

$ws.Range("H138").Value = 12479.777
$ws.Range("I138").Value = 13762.625
$ws.Range("J138").Value = 2217
$ws.Range("K138").Value = 41287.875
$ws.Range("L138").Value = 6651
$ws.Range("M138").Value = -36147.875
$ws.Range("N138").Value = -16931

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 2015.0834
$ws.Range("I5").Value = 256.8
$ws.Range("J5").Value = 3271
$ws.Range("K5").Value = 256.8
$ws.Range("L5").Value = 3271
$ws.Range("M5").Value = -144.8
$ws.Range("N5").Value = -3495

$ws.Range("H46").Value = 8000
$ws.Range("J46").Value = 11000
$ws.Range("L46").Value = 11000
$ws.Range("N46").Value = -11638

$ws.Range("H101").Value = 317285.88
$ws.Range("J101").Value = 317285.88
$ws.Range("L101").Value = 317285.88
$ws.Range("N101").Value = -323775.88

$ws.Range("H132").Value = 2044393
$ws.Range("I132").Value = 2581402
$ws.Range("K132").Value = 7744206
$ws.Range("M132").Value = -7741676

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 2015.0834
$ws.Range("I4").Value = 256.8
$ws.Range("J4").Value = 3271
$ws.Range("K4").Value = 256.8
$ws.Range("L4").Value = 3271
$ws.Range("M4").Value = -141.8
$ws.Range("N4").Value = -3501

$ws.Range("H86").Value = 3729.2222
$ws.Range("I86").Value = 3916.5
$ws.Range("K86").Value = 3916.5
$ws.Range("M86").Value = -2793.5

$ws.Range("H89").Value = 3729.2222
$ws.Range("I89").Value = 3916.5
$ws.Range("K89").Value = 19582.5
$ws.Range("M89").Value = -13966.5

$ws.Range("H134").Value = 4467
$ws.Range("I134").Value = 2140.7104
$ws.Range("K134").Value = 6422.1312
$ws.Range("M134").Value = -3887.1312

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1512.1666
$ws.Range("I7").Value = 45.666668
$ws.Range("K7").Value = 45.666668
$ws.Range("M7").Value = 67.333332

$ws.Range("H22").Value = 1430.2916
$ws.Range("I22").Value = 177.90909
$ws.Range("J22").Value = 2490
$ws.Range("K22").Value = 177.90909
$ws.Range("L22").Value = 2490
$ws.Range("M22").Value = 172.09091
$ws.Range("N22").Value = -3190

$ws.Range("H104").Value = 68000
$ws.Range("J104").Value = 68000
$ws.Range("L104").Value = 68000
$ws.Range("N104").Value = -73242

$ws.Range("H132").Value = 1706.814
$ws.Range("I132").Value = 1167.8948
$ws.Range("K132").Value = 3503.6844
$ws.Range("M132").Value = -973.6844000000001

$ws.Range("H134").Value = 5113.1665
$ws.Range("I134").Value = 3919.75
$ws.Range("K134").Value = 11759.25
$ws.Range("M134").Value = -9224.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 66666812
$ws.Range("I11").Value = 99.55556
$ws.Range("J11").Value = 166666880
$ws.Range("K11").Value = 298.66668
$ws.Range("L11").Value = 500000640
$ws.Range("M11").Value = -158.66668
$ws.Range("N11").Value = -500000920

$ws.Range("H99").Value = 1735.6
$ws.Range("I99").Value = 1428.5
$ws.Range("J99").Value = 2964
$ws.Range("K99").Value = 4285.5
$ws.Range("L99").Value = 8892
$ws.Range("M99").Value = -2039.5
$ws.Range("N99").Value = -13384

$ws.Range("H136").Value = 10812.125
$ws.Range("I136").Value = 10000
$ws.Range("K136").Value = 30000
$ws.Range("M136").Value = -24900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 8964461
$ws.Range("I11").Value = 8964461
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 8964461
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -8964322
$ws.Range("N11").ClearContents()

$ws.Range("H122").Value = 4176.5
$ws.Range("I122").Value = 4011.8
$ws.Range("K122").Value = 12035.4
$ws.Range("M122").Value = -9585.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1300.8889
$ws.Range("I9").Value = 213.5
$ws.Range("K9").Value = 213.5
$ws.Range("M9").Value = 10.5

$ws.Range("H20").Value = 46666.668
$ws.Range("I20").Value = 17000
$ws.Range("K20").Value = 17000
$ws.Range("M20").Value = -16774

$ws.Range("H122").Value = 2608.111
$ws.Range("J122").Value = 3140.8
$ws.Range("L122").Value = 9422.400000000001
$ws.Range("N122").Value = -14322.4

$ws.Range("H132").Value = 2755.9807
$ws.Range("I132").Value = 2226.1714
$ws.Range("K132").Value = 6678.514200000001
$ws.Range("M132").Value = -4148.514200000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5259.8
$ws.Range("J62").Value = 5949.5
$ws.Range("L62").Value = 5949.5
$ws.Range("N62").Value = -7197.5

$ws.Range("H65").Value = 5259.8
$ws.Range("J65").Value = 5949.5
$ws.Range("L65").Value = 29747.5
$ws.Range("N65").Value = -35987.5

$ws.Range("H107").Value = 2787.0476
$ws.Range("I107").Value = 1220
$ws.Range("J107").Value = 9447
$ws.Range("K107").Value = 3660
$ws.Range("L107").Value = 28341
$ws.Range("M107").Value = -1740
$ws.Range("N107").Value = -32181
